$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix data errors ---
# Row 4: ID was text "a" -> should be numeric 3
$ws.Cells.Item(4, 1).Value2 = 3

# Row 5: NOMBRE was "Adri" -> should be "sdfs"
$ws.Cells.Item(5, 2).Value2 = "sdfs"

# Row 6: Email column had wrong text "alan@gmail.com" -> should read "456456@"
# (hyperlink target itself is left untouched by the source diff)
# Give it the "ErrorColumns" hyperlink look (underlined hyperlink style)
$ws.Cells.Item(6, 4).Style = "Hipervínculo"
$ws.Cells.Item(6, 4).Value2 = "456456@"

# Row 8: ID was text "gh" -> should be numeric 7
$ws.Cells.Item(8, 1).Value2 = 7

# Row 8: Phone 43534 -> 4353
$ws.Cells.Item(8, 5).Value2 = 4353

# Row 11: Phone was empty -> should be 324
$ws.Cells.Item(11, 5).Value2 = 324

# Row 12: ID 11 -> 33
$ws.Cells.Item(12, 1).Value2 = 33

# --- Selection bookkeeping (matches saved workbook view state) ---
$ws.Range("I15").Select()
